$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: "Change Management Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Change Management Overview")

# Project name: AI/ML -> Cloud Infrastructure Migration
$ws1.Range("B6").Value = "Enterprise Cloud Infrastructure Migration"

# Objectives: drop the AI/ML framing in favor of generic "IT"
$ws1.Range("A15").Value = "1. Achieve 95% user adoption of new IT systems within 6 months of go-live"
$ws1.Range("A17").Value = "3. Build organizational capability and confidence in IT technologies"
$ws1.Range("A20").Value = "6. Create positive stakeholder sentiment and enthusiasm for IT transformation"

# Preserve the pre-existing blank row (4) so the save round-trip doesn't
# silently drop it, and restore the blank spacer rows (13 and 21) that
# separate the sections below - all without shifting any row numbers.
$ws1.Rows.Item(4).Group()
$ws1.Rows.Item(4).Ungroup()
$ws1.Rows.Item(13).Group()
$ws1.Rows.Item(13).Ungroup()
$ws1.Rows.Item(21).Group()
$ws1.Rows.Item(21).Ungroup()

# ---------------------------------------------------------------------------
# Sheet: "Change Impact Assessment"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Change Impact Assessment")

$ws2.Range("A4").Value = "IT Managers"
$ws2.Range("G4").Value = "IT automation"
$ws2.Range("A5").Value = "System Administrators"

# Restore the blank spacer row (2) between the title and the header row.
$ws2.Rows.Item(2).Group()
$ws2.Rows.Item(2).Ungroup()

# ---------------------------------------------------------------------------
# Sheet: "Change Activities"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Change Activities")

# Restore the blank spacer row (2) between the title and the header row.
$ws3.Rows.Item(2).Group()
$ws3.Rows.Item(2).Ungroup()
